$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before BB (old BB/BC shift right to become BC/BD)
$ws.Columns("BB:BB").Insert()

# New header for the inserted timestamp column
$ws.Range("BB1").Value = "2026-01-30 03:26:01"

# Populate the new BB column for each data row: it mirrors the latest
# price snapshot that was already recorded in column BA for that row.
# Rows whose BA cell has no price (inlineStr/blank) are left untouched.
for ($r = 2; $r -le 206; $r++) {
    $baCell = $ws.Cells.Item($r, 53)
    if ($baCell.Value2 -ne $null -and $baCell.Value2 -ne "") {
        $ws.Cells.Item($r, 54).Value = $baCell.Value2
    }
}
